# The deck ships two theme parts:
#   theme1.xml -> "Office" colour scheme (used only by the Notes Master)
#   theme2.xml -> "Red Violet"/"Integral" colour scheme (used by the one
#                 Slide Master that backs every slide layout/slide)
#
# The authored change swaps the two themes' contents so that the design
# actually applied to the slides (theme2.xml) becomes the standard
# "Office" palette, while the "Integral"/Red-Violet palette moves into
# theme1.xml (the Notes Master's theme).
#
# This host exposes the live (slide-facing) theme's colour scheme via
# Slide.ThemeColorScheme / Slide.NotesPage.ThemeColorScheme, in the
# fixed MsoThemeColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# Writing RGB there updates the backing theme XML (theme2.xml) in place,
# which is how we recreate the "Office" colours that PowerPoint's Design
# gallery would have swapped in.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Target palette = the standard Office theme colours (the content that
# used to live in theme1.xml, now authored into theme2.xml).
$cs.Item(1).RGB  = 0        # dk1      000000
$cs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388  # dk2      44546A
$cs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501  # accent2  ED7D31
$cs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$cs.Item(8).RGB  = 49407    # accent4  FFC000
$cs.Item(9).RGB  = 12874308 # accent5  4472C4
$cs.Item(10).RGB = 4697456  # accent6  70AD47
$cs.Item(11).RGB = 12673797 # hlink    0563C1
$cs.Item(12).RGB = 7491477  # folHlink 954F72
